# room.xlsx update: renumber rooms to a contiguous sequence per floor,
# fix a couple of vt_no (column E) values that were swapped, move the
# active selection, and set the page to A4 portrait for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (room_no): renumber rooms 306/308/310/312 -> 305/306/307/308, etc. ---
$ws.Range("A6").Value  = 305
$ws.Range("A7").Value  = 306
$ws.Range("A8").Value  = 307
$ws.Range("A9").Value  = 308

$ws.Range("A12").Value = 403
$ws.Range("A13").Value = 404
$ws.Range("A14").Value = 405
$ws.Range("A15").Value = 406
$ws.Range("A16").Value = 407
$ws.Range("A17").Value = 408

$ws.Range("A20").Value = 503
$ws.Range("A21").Value = 504
$ws.Range("A22").Value = 505
$ws.Range("A23").Value = 506

$ws.Range("A26").Value = 603
$ws.Range("A27").Value = 604
$ws.Range("A28").Value = 605

# --- Column E (vt_no) corrections ---
$ws.Range("E25").Value = 2
$ws.Range("E26").Value = 1
$ws.Range("E31").Value = 2
$ws.Range("E32").Value = 1

# --- Move the active selection to C13 ---
$ws.Range("C13").Select()

# --- Page setup: A4, portrait ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
